$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.024.62'
$ws.Range('E2').Value = '  +3.05%  '
$ws.Range('D3').Value = '2.454.73'
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.33'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.67'
$ws.Range('E6').Value = '  +3.36%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').Value = '2.453.99'
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('E10').Value = '  +2.82%  '
$ws.Range('E11').Value = '  +2.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').Value = '  +1.78%  '
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('E14').Value = '  +7.41%  '
$ws.Range('E15').Value = '  +5.67%  '
$ws.Range('D16').Value = '2.898.67'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D17').Value = '62.797.98'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').Value = '2.451.59'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('E20').Value = '  +4.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '330.94'
$ws.Range('E21').Value = '  +2.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.14'
$ws.Range('E23').Value = '  +11.51%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.42'
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('E26').Value = '  +24.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '650.95'
$ws.Range('E27').Value = '  +10.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.65'
$ws.Range('E28').Value = '  +5.16%  '
$ws.Range('E29').Value = '  +5.76%  '
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.20'
$ws.Range('E31').Value = '  +2.39%  '
$ws.Range('E32').Value = '  +5.63%  '
$ws.Range('E33').Value = '  +3.82%  '
$ws.Range('E34').Value = '  +4.17%  '
$ws.Range('D35').Value = '0.0₆0419'
$ws.Range('E35').Value = '  +48.49%  '
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.79'
$ws.Range('E38').Value = '  +3.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.56'
$ws.Range('E39').Value = '  +5.98%  '
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '152.66'
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.84'
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.73'
$ws.Range('E43').Value = '  +9.98%  '
$ws.Range('E44').Value = '  +4.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.63'
$ws.Range('E45').Value = '  +2.17%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.00'
$ws.Range('E47').Value = '  +27.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '146.47'
$ws.Range('E48').Value = '  +3.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.75'
$ws.Range('E50').Value = '  +4.99%  '
$ws.Range('E51').Value = '  +2.52%  '
